$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the revision date banner
$ws.Range("A1").Value = "dSIB  Revised: January 26, 2017"

# Update the "Bill Of Materials" title/timestamp line
$ws.Range("A8").Value = "Bill Of Materials       January 26, 2017"

# Update D3 part (row 25): swap Diodes Inc. ZHCS750/SOT for ST Micro BAR46FILM,
# with a new vendor part number and a 100V voltage rating.
$ws.Range("D25").Value = "BAR46FILM"
$ws.Range("E25").Value = "ST Micro"
$ws.Range("F25").Value = "BAR46FILM"
$ws.Range("H25").Value = "497-12128-1-ND"
$ws.Range("K25").Value = "100V"

# Update the active cell selection
$ws.Range("A9").Select()
